$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header (shared string) values
$ws.Range("B1").Value = "TISG"
$ws.Range("D1").Value = "buy"
$ws.Range("E1").Value = "MYDIR"
$ws.Range("F1").Value = "need_to_buy"

# New data values for rows 2-15 (A: date, B:F numeric)
$data = @(
    @(45964, 7938.18802635075, 7065.41076563385, 9004, 3526.04817118358, 66.1441223673929),
    @(45965, 8368.56154267988, 7946.8043023883, 9004, 4492.54431047616, 143.139525536019),
    @(45966, 8591.6560040422, 8232.3434855197, 9004, 4638.28571216266, 161.109549903432),
    @(45967, 8575.06798513402, 8378.424989253, 9004, 4611.19407830248, 166.067461148145),
    @(45968, 8399.71252085342, 7554.20967854608, 9004, 4527.70949005973, 128.246632025242),
    @(45969, 3116.75094940757, 5339.01661445299, 9004, 4304.68543734191, 26.6542521581206),
    @(45970, 2982.76489981526, 4951.12846317757, 9004, 4297.42657850312, 10.1897934033621),
    @(45971, 9273.47168962819, 8225.3796329613, 9004, 5243.41242511772, 186.033002419959),
    @(45972, 9273.47299076555, 8471.19158803462, 9004, 5243.41618746352, 196.275323979089),
    @(45973, 9273.47299076555, 8629.97273475286, 9004, 5243.41618746352, 202.891205092349),
    @(45974, 9273.47299076555, 8734.5131946171, 9004, 5243.41618746352, 207.247057586692),
    @(45975, 9273.47299076555, 7827.90534151468, 9004, 5243.41618746352, 169.471730374091),
    @(45976, 3551.20396235217, 5021.47838587834, 9004, 4928.52850584403, 39.4169538217657),
    @(45977, 3412.56435984501, 4786.18433598882, 9004, 4921.05669262868, 29.3017095257293)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
